# Fix todo list to reflect last commit...
#
# - Row 23: replace "Guide and Channels from OSD need the Back button at top
#   left." with "Use full guide rather than the abbreviated OSD guide.
#   Remove OSD channel button." and mark it "Done" (B23) instead of having a
#   Comments entry in C23.
# - Row 27 (Comments, C27): update the note about merging Guide/Channels.
# - Scroll/selection state changes to topLeftCell A16 / selection B27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 23: update feature text, clear the old Comments cell, and set Status
$ws.Range("A23").Value = "Use full guide rather than the abbreviated OSD guide. Remove OSD channel button."
$ws.Range("C23").Clear()
$ws.Range("B23").Value = "Done"

# Row 27: update Comments text
$ws.Range("C27").Value = "Can we just have the Guide and get rid of Channels? For the moment, just make guide the first widget in the line."

# Update the view's scroll position and selection
$ws.Range("B27").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
